$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on Price (D) cells before assigning, then reset the
# cell style back to "Normal" so the written value stays an exact string (no
# Excel auto-numeric-coercion / trailing-zero loss) while no residual style index
# is left attached to the cell (matches the source, which has no s= attribute).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.051.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.70%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.37%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.70%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3012"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06900"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.914.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07310"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.46%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.53%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.111"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.20%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6835"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.50%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.045.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008087"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.169.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.888"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "184.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +34.91%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.098"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.82%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.399"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.86%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.954"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.408"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.366"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.16%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08979"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.069"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05263"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7523"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.150"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.77%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.671"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01941"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +18.17%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.743"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.196"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9394"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4390"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.56%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.897"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.20%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.827"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1341"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.14%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05867"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.658"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.66%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3907"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.402"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.83%  "
